$d = $word.ActiveDocument

# Helper: insert a new paragraph right after paragraph index $afterIdx,
# set its text to $text, and return the index of the newly created paragraph
# (so subsequent calls can chain off of it).
function Insert-Line($afterIdx, $text) {
    $d.Paragraphs($afterIdx).Range.InsertParagraphAfter() | Out-Null
    $d.Paragraphs($afterIdx + 1).Range.Text = $text
    return $afterIdx + 1
}

# 1) <title> text change
$d.Content.Find.Execute("Mi Primera Página", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "COMBEE - Landing Page", 2) | Out-Null

# 2) body { ... } block: "margin: 20px;" -> "background-color: #f0f0f0;" plus
#    several new declarations after it.
$d.Content.Find.Execute("            margin: 20px;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "            background-color: #f0f0f0;", 2) | Out-Null

$idx = 8
$idx = Insert-Line $idx "            margin: 0;"
$idx = Insert-Line $idx "            padding: 0;"
$idx = Insert-Line $idx "            display: flex;"
$idx = Insert-Line $idx "            justify-content: center;"
$idx = Insert-Line $idx "            align-items: center;"
$idx = Insert-Line $idx "            height: 100vh;"
$idx = Insert-Line $idx "            text-align: center;"
# $idx now points at the last inserted line; paragraph $idx+1 is the original "        }"

# 3) New ".container { ... }" rule block, inserted right before "        h1 {"
$idx = $idx + 1  # now at "        }" (closing body rule)
$idx = Insert-Line $idx "        .container {"
$idx = Insert-Line $idx "            background-color: #ffffff;"
$idx = Insert-Line $idx "            padding: 50px;"
$idx = Insert-Line $idx "            border-radius: 10px;"
$idx = Insert-Line $idx "            box-shadow: 0 4px 8px rgba(0, 0, 0, 0.1);"
$idx = Insert-Line $idx "        }"
# paragraph $idx+1 is the original "        h1 {"

# 4) h1 { color: blue; } -> color: #FFD700; plus font-size: 48px;
$d.Content.Find.Execute("            color: blue;", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "            color: #FFD700; /* Color amarillo */", 2) | Out-Null

# Locate a paragraph index by (unique) substring - used to avoid re-deriving indices by hand.
function Get-ParaIndex($searchText) {
    for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
        if ($d.Paragraphs($i).Range.Text -like "*$searchText*") {
            return $i
        }
    }
    return -1
}

$idx = Get-ParaIndex "color: #FFD700"
$idx = Insert-Line $idx "            font-size: 48px;"

# 5) p { font-size: 18px; } -> add color: #333333; right before it
$idx = Get-ParaIndex "font-size: 18px;"
$idx = Insert-Line ($idx - 1) "            color: #333333;"

# 6) <body> content: replace <h1> and <p> lines with a wrapping <div class="container"> block
$idx = Get-ParaIndex "Esta es mi primera p"
$d.Paragraphs($idx).Range.Delete() | Out-Null

$idx = Get-ParaIndex "Hola, mundo"
$d.Paragraphs($idx).Range.Text = "    <div class=`"container`">"
$idx = Insert-Line $idx "        <h1>COMBEE</h1>"
$idx = Insert-Line $idx "        <p>Bienvenido a la landing page de COMBEE. Estamos trabajando para ofrecerte la mejor experiencia de autos compartidos.</p>"
$idx = Insert-Line $idx "    </div>"
